$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counts for the existing rows
$ws.Range("B2").Value = 79224
$ws.Range("B3").Value = 2750
$ws.Range("B4").Value = 542

# Remove row 5 entirely (monzodiorite monzogabbro / 140)
$ws.Rows.Item(5).Delete()
